$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3768554877997065
$ws.Range("J2").Value = 0.3768554877997065
$ws.Range("M2").Value = 8.252454666666667
$ws.Range("N2").Value = 24.757364
$ws.Range("O2").Value = 0.05349680956196952
$ws.Range("P2").Value = 0.05349680956196953
$ws.Range("Q2").Value = 1.222295818044
$ws.Range("R2").Value = 11.000662362396
$ws.Range("S2").Value = 0.02016056626320403
$ws.Range("T2").Value = 0.02016056626320403

$ws.Range("I3").Value = 0.3768554877997065
$ws.Range("J3").Value = 0.3768554877997065
$ws.Range("O3").Value = 0.5638948237978928
$ws.Range("P3").Value = 0.5638948237978929
$ws.Range("S3").Value = 0.2125068588900844
$ws.Range("T3").Value = 0.2125068588900845

$ws.Range("I4").Value = 0.3768554877997065
$ws.Range("J4").Value = 0.3768554877997065
$ws.Range("M4").Value = 57.81408433333333
$ws.Range("N4").Value = 173.442253
$ws.Range("O4").Value = 0.3747817085348802
$ws.Range("P4").Value = 0.3747817085348802
$ws.Range("Q4").Value = 8.563017472862999
$ws.Range("R4").Value = 77.067157255767
$ws.Range("S4").Value = 0.1412385435883197
$ws.Range("T4").Value = 0.1412385435883197

$ws.Range("I5").Value = 0.3768554877997065
$ws.Range("J5").Value = 0.3768554877997065
$ws.Range("M5").Value = 1.207345666666667
$ws.Range("N5").Value = 3.622037
$ws.Range("O5").Value = 0.007826658105257385
$ws.Range("P5").Value = 0.007826658105257386
$ws.Range("Q5").Value = 0.178823588727
$ws.Range("R5").Value = 1.609412298543
$ws.Range("S5").Value = 0.002949519058098299
$ws.Range("T5").Value = 0.002949519058098299

$ws.Range("G6").Value = 0.2449103333333333
$ws.Range("H6").Value = 0.734731
$ws.Range("I6").Value = 0.6231445122002934
$ws.Range("J6").Value = 0.6231445122002934
$ws.Range("M6").Value = 8.252454666666667
$ws.Range("N6").Value = 24.757364
$ws.Range("O6").Value = 0.05349680956196952
$ws.Range("P6").Value = 0.05349680956196953
$ws.Range("Q6").Value = 2.021111423231556
$ws.Range("R6").Value = 18.190002809084
$ws.Range("S6").Value = 0.03333624329876549
$ws.Range("T6").Value = 0.03333624329876549

$ws.Range("G7").Value = 0.2449103333333333
$ws.Range("H7").Value = 0.734731
$ws.Range("I7").Value = 0.6231445122002934
$ws.Range("J7").Value = 0.6231445122002934
$ws.Range("O7").Value = 0.5638948237978928
$ws.Range("P7").Value = 0.5638948237978929
$ws.Range("Q7").Value = 21.30396708160456
$ws.Range("R7").Value = 191.735703734441
$ws.Range("S7").Value = 0.3513879649078083
$ws.Range("T7").Value = 0.3513879649078084

$ws.Range("G8").Value = 0.2449103333333333
$ws.Range("H8").Value = 0.734731
$ws.Range("I8").Value = 0.6231445122002934
$ws.Range("J8").Value = 0.6231445122002934
$ws.Range("M8").Value = 57.81408433333333
$ws.Range("N8").Value = 173.442253
$ws.Range("O8").Value = 0.3747817085348802
$ws.Range("P8").Value = 0.3747817085348802
$ws.Range("Q8").Value = 14.15926666543811
$ws.Range("R8").Value = 127.433399988943
$ws.Range("S8").Value = 0.2335431649465605
$ws.Range("T8").Value = 0.2335431649465605

$ws.Range("G9").Value = 0.2449103333333333
$ws.Range("H9").Value = 0.734731
$ws.Range("I9").Value = 0.6231445122002934
$ws.Range("J9").Value = 0.6231445122002934
$ws.Range("M9").Value = 1.207345666666667
$ws.Range("N9").Value = 3.622037
$ws.Range("O9").Value = 0.007826658105257385
$ws.Range("P9").Value = 0.007826658105257386
$ws.Range("Q9").Value = 0.2956914296718889
$ws.Range("R9").Value = 2.661222867047
$ws.Range("S9").Value = 0.004877139047159086
$ws.Range("T9").Value = 0.004877139047159086
